$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Cells.Item(11, 8).Value = 82.78570999999999
$ws.Cells.Item(11, 9).Value = 82.78570999999999
$ws.Cells.Item(11, 11).Value = 82.78570999999999
$ws.Cells.Item(11, 13).Value = 57.21429000000001
# Row 62
$ws.Cells.Item(62, 8).Value = 2999
$ws.Cells.Item(62, 9).Value = 2999
$ws.Cells.Item(62, 11).Value = 2999
$ws.Cells.Item(62, 13).Value = -2375
# Row 65
$ws.Cells.Item(65, 8).Value = 2999
$ws.Cells.Item(65, 9).Value = 2999
$ws.Cells.Item(65, 11).Value = 14995
$ws.Cells.Item(65, 13).Value = -11875
# Row 86
$ws.Cells.Item(86, 8).Value = 3025.476
$ws.Cells.Item(86, 10).Value = 5204.5713
$ws.Cells.Item(86, 12).Value = 5204.5713
$ws.Cells.Item(86, 14).Value = -7450.5713
# Row 89
$ws.Cells.Item(89, 8).Value = 3025.476
$ws.Cells.Item(89, 10).Value = 5204.5713
$ws.Cells.Item(89, 12).Value = 26022.8565
$ws.Cells.Item(89, 14).Value = -37254.85649999999
# Row 92
$ws.Cells.Item(92, 8).Value = 1869.7778
$ws.Cells.Item(92, 9).Value = 2007.2142
$ws.Cells.Item(92, 10).Value = 1388.75
$ws.Cells.Item(92, 11).Value = 2007.2142
$ws.Cells.Item(92, 12).Value = 1388.75
$ws.Cells.Item(92, 13).Value = -759.2141999999999
$ws.Cells.Item(92, 14).Value = -3884.75
# Row 98
$ws.Cells.Item(98, 8).Value = 1795.875
$ws.Cells.Item(98, 9).Value = 1795.875
$ws.Cells.Item(98, 11).Value = 1795.875
$ws.Cells.Item(98, 13).Value = -297.875
# Row 103
$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 0
$ws.Cells.Item(103, 13).ClearContents()
$ws.Cells.Item(103, 14).ClearContents()
# Row 107
$ws.Cells.Item(107, 8).Value = 1530.5625
$ws.Cells.Item(107, 9).Value = 1999.7273
$ws.Cells.Item(107, 10).Value = 498.4
$ws.Cells.Item(107, 11).Value = 1999.7273
$ws.Cells.Item(107, 12).Value = 498.4
$ws.Cells.Item(107, 13).Value = -79.72730000000001
$ws.Cells.Item(107, 14).Value = -4338.4
# Row 122
$ws.Cells.Item(122, 8).Value = 1795.875
$ws.Cells.Item(122, 9).Value = 1795.875
$ws.Cells.Item(122, 11).Value = 5387.625
$ws.Cells.Item(122, 13).Value = -2937.625
# Row 130
$ws.Cells.Item(130, 8).Value = 24998.166
$ws.Cells.Item(130, 10).Value = 24998.166
$ws.Cells.Item(130, 12).Value = 24998.166
$ws.Cells.Item(130, 14).Value = -35038.166
# Row 132
$ws.Cells.Item(132, 8).Value = 15630.359
$ws.Cells.Item(132, 9).Value = 2433.9062
$ws.Cells.Item(132, 11).Value = 7301.7186
$ws.Cells.Item(132, 13).Value = -4771.7186
# Row 138
$ws.Cells.Item(138, 8).Value = 2760.311
$ws.Cells.Item(138, 9).Value = 1471.0834
$ws.Cells.Item(138, 11).Value = 4413.2502
$ws.Cells.Item(138, 13).Value = 726.7497999999996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 33250
$ws.Cells.Item(61, 9).Value = 33250
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 33250
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 14).ClearContents()
# Row 74
$ws.Cells.Item(74, 8).Value = 1712.7693
$ws.Cells.Item(74, 9).Value = 1771
$ws.Cells.Item(74, 11).Value = 1771
$ws.Cells.Item(74, 13).Value = -897
# Row 77
$ws.Cells.Item(77, 8).Value = 1712.7693
$ws.Cells.Item(77, 9).Value = 1771
$ws.Cells.Item(77, 11).Value = 8855
$ws.Cells.Item(77, 13).Value = -4487
# Row 97
$ws.Cells.Item(97, 8).Value = 2597.6667
$ws.Cells.Item(97, 9).Value = 1165.0714
$ws.Cells.Item(97, 10).Value = 7611.75
$ws.Cells.Item(97, 11).Value = 1165.0714
$ws.Cells.Item(97, 12).Value = 7611.75
$ws.Cells.Item(97, 13).Value = -669.0714
$ws.Cells.Item(97, 14).Value = -8603.75
# Row 132
$ws.Cells.Item(132, 8).Value = 3806.3572
$ws.Cells.Item(132, 9).Value = 3806.3572
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 11419.0716
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 14).ClearContents()
# Row 136
$ws.Cells.Item(136, 8).Value = 33250
$ws.Cells.Item(136, 9).Value = 33250
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 99750
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 1921.0741
$ws.Cells.Item(94, 9).Value = 845.7895
$ws.Cells.Item(94, 11).Value = 845.7895
$ws.Cells.Item(94, 13).Value = -394.7895
# Row 105
$ws.Cells.Item(105, 8).Value = 2279.4285
$ws.Cells.Item(105, 9).Value = 797.93335
$ws.Cells.Item(105, 11).Value = 797.93335
$ws.Cells.Item(105, 13).Value = 949.06665
# Row 107
$ws.Cells.Item(107, 8).Value = 7629.32
$ws.Cells.Item(107, 9).Value = 7286.15
$ws.Cells.Item(107, 10).Value = 9002
$ws.Cells.Item(107, 11).Value = 7286.15
$ws.Cells.Item(107, 12).Value = 9002
$ws.Cells.Item(107, 13).Value = -5366.15
$ws.Cells.Item(107, 14).Value = -12842
# Row 134
$ws.Cells.Item(134, 8).Value = 4347.6665
$ws.Cells.Item(134, 9).Value = 3932.25
$ws.Cells.Item(134, 11).Value = 11796.75
$ws.Cells.Item(134, 13).Value = -9261.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 1363.0638
$ws.Cells.Item(31, 9).Value = 1293.6444
$ws.Cells.Item(31, 10).Value = 2925
$ws.Cells.Item(31, 11).Value = 1293.6444
$ws.Cells.Item(31, 12).Value = 2925
$ws.Cells.Item(31, 13).Value = -998.6443999999999
$ws.Cells.Item(31, 14).Value = -3515
# Row 34
$ws.Cells.Item(34, 8).Value = 1363.0638
$ws.Cells.Item(34, 9).Value = 1293.6444
$ws.Cells.Item(34, 10).Value = 2925
$ws.Cells.Item(34, 11).Value = 1293.6444
$ws.Cells.Item(34, 12).Value = 2925
$ws.Cells.Item(34, 13).Value = -1091.6444
$ws.Cells.Item(34, 14).Value = -3329
# Row 105
$ws.Cells.Item(105, 8).Value = 950.3570999999999
$ws.Cells.Item(105, 9).Value = 928
$ws.Cells.Item(105, 11).Value = 928
$ws.Cells.Item(105, 13).Value = 819
# Row 132
$ws.Cells.Item(132, 8).Value = 2384.5
$ws.Cells.Item(132, 9).Value = 1716.1177
$ws.Cells.Item(132, 10).Value = 4657
$ws.Cells.Item(132, 11).Value = 5148.3531
$ws.Cells.Item(132, 12).Value = 13971
$ws.Cells.Item(132, 13).Value = -2618.3531
$ws.Cells.Item(132, 14).Value = -19031
# Row 134
$ws.Cells.Item(134, 8).Value = 2500
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 2500
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 13).ClearContents()
$ws.Cells.Item(134, 14).Value = -12570
# Row 141
$ws.Cells.Item(141, 8).Value = 47875
$ws.Cells.Item(141, 9).Value = 21000
$ws.Cells.Item(141, 11).Value = 21000
$ws.Cells.Item(141, 13).Value = -15820

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Cells.Item(12, 8).Value = 8310.532999999999
$ws.Cells.Item(12, 9).Value = 1475
$ws.Cells.Item(12, 10).Value = 10796.182
$ws.Cells.Item(12, 11).Value = 4425
$ws.Cells.Item(12, 12).Value = 32388.546
$ws.Cells.Item(12, 13).Value = -4252
$ws.Cells.Item(12, 14).Value = -32734.546
# Row 18
$ws.Cells.Item(18, 8).Value = 321.9
$ws.Cells.Item(18, 9).Value = 152.71428
$ws.Cells.Item(18, 10).Value = 716.6667
$ws.Cells.Item(18, 11).Value = 458.14284
$ws.Cells.Item(18, 12).Value = 2150.0001
$ws.Cells.Item(18, 13).Value = -289.14284
$ws.Cells.Item(18, 14).Value = -2488.0001
# Row 113
$ws.Cells.Item(113, 8).Value = 1717.1538
$ws.Cells.Item(113, 9).Value = 1957.8334
$ws.Cells.Item(113, 10).Value = 1510.8572
$ws.Cells.Item(113, 11).Value = 5873.5002
$ws.Cells.Item(113, 12).Value = 4532.571599999999
$ws.Cells.Item(113, 13).Value = -3703.5002
$ws.Cells.Item(113, 14).Value = -8872.571599999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 5975.8125
$ws.Cells.Item(80, 10).Value = 9005.571
$ws.Cells.Item(80, 12).Value = 9005.571
$ws.Cells.Item(80, 14).Value = -11001.571
# Row 83
$ws.Cells.Item(83, 8).Value = 5975.8125
$ws.Cells.Item(83, 10).Value = 9005.571
$ws.Cells.Item(83, 12).Value = 45027.855
$ws.Cells.Item(83, 14).Value = -55011.855
# Row 113
$ws.Cells.Item(113, 8).Value = 7941.357
$ws.Cells.Item(113, 9).Value = 6169.7144
$ws.Cells.Item(113, 10).Value = 9713
$ws.Cells.Item(113, 11).Value = 6169.7144
$ws.Cells.Item(113, 12).Value = 9713
$ws.Cells.Item(113, 13).Value = -3999.7144
$ws.Cells.Item(113, 14).Value = -14053
# Row 132
$ws.Cells.Item(132, 8).Value = 7116.6343
$ws.Cells.Item(132, 9).Value = 6232.8237
$ws.Cells.Item(132, 11).Value = 18698.4711
$ws.Cells.Item(132, 13).Value = -16168.4711

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 43
$ws.Cells.Item(43, 8).Value = 38937.5
$ws.Cells.Item(43, 9).Value = 39333.332
$ws.Cells.Item(43, 10).Value = 38846.152
$ws.Cells.Item(43, 11).Value = 39333.332
$ws.Cells.Item(43, 12).Value = 38846.152
$ws.Cells.Item(43, 13).Value = -39140.332
$ws.Cells.Item(43, 14).Value = -39232.152
# Row 100
$ws.Cells.Item(100, 8).Value = 6717
$ws.Cells.Item(100, 9).Value = 4945.3335
$ws.Cells.Item(100, 10).Value = 7780
$ws.Cells.Item(100, 11).Value = 4945.3335
$ws.Cells.Item(100, 12).Value = 7780
$ws.Cells.Item(100, 13).Value = -4404.3335
$ws.Cells.Item(100, 14).Value = -8862
# Row 127
$ws.Cells.Item(127, 8).Value = 54799.5
$ws.Cells.Item(127, 10).Value = 54799.5
$ws.Cells.Item(127, 12).Value = 54799.5
$ws.Cells.Item(127, 14).Value = -64719.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Cells.Item(51, 8).Value = 25943.555
$ws.Cells.Item(51, 10).Value = 33498.6
$ws.Cells.Item(51, 12).Value = 33498.6
$ws.Cells.Item(51, 14).Value = -34518.6
# Row 61
$ws.Cells.Item(61, 8).Value = 44799
$ws.Cells.Item(61, 10).Value = 44799
$ws.Cells.Item(61, 12).Value = 44799
$ws.Cells.Item(61, 14).Value = -45383
# Row 113
$ws.Cells.Item(113, 8).Value = 412.58334
$ws.Cells.Item(113, 9).Value = 429
$ws.Cells.Item(113, 10).Value = 363.33334
$ws.Cells.Item(113, 11).Value = 1287
$ws.Cells.Item(113, 12).Value = 1090.00002
$ws.Cells.Item(113, 13).Value = 883
$ws.Cells.Item(113, 14).Value = -5430.000019999999
# Row 136
$ws.Cells.Item(136, 8).Value = 1130.4615
$ws.Cells.Item(136, 9).Value = 1099.6666
$ws.Cells.Item(136, 10).Value = 1500
$ws.Cells.Item(136, 11).Value = 3298.9998
$ws.Cells.Item(136, 12).Value = 4500
$ws.Cells.Item(136, 13).Value = -748.9998000000001
$ws.Cells.Item(136, 14).Value = -9600
